# Applies odds updates to Sheet1 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("BD3").Value = 151

# Row 5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25

# Row 6
$ws.Range("G6").Value = 2.15
$ws.Range("I6").Value = 3.7
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 19
$ws.Range("AE6").Value = 17
$ws.Range("AW6").Value = 5.5

# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
